$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp update
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 12:30"

# India (row 7)
$ws.Range("B7").Value = 650431
$ws.Range("C7").Value = 542
$ws.Range("E7").Value = 237351

# Rows 30/31: Indonesia overtakes Belgica in ranking -> swap country names,
# row 30 gets Indonesia's new totals, row 31 gets Belgica's (previous) totals.
$ws.Range("A30").Value = "Indonesia"
$ws.Range("B30").Value = 62142
$ws.Range("C30").Value = 1447
$ws.Range("D30").Value = 28219
$ws.Range("E30").Value = 30834
$ws.Range("G30").Value = 53
$ws.Range("H30").Value = 3089

$ws.Range("A31").Value = "Belgica"
$ws.Range("B31").Value = 61727
$ws.Range("D31").Value = 17073
$ws.Range("E31").Value = 34889
$ws.Range("H31").Value = 9765

# Filipinas (row 42)
$ws.Range("B42").Value = 41830
$ws.Range("C42").Value = 1494
$ws.Range("D42").Value = 11453
$ws.Range("E42").Value = 29087
$ws.Range("G42").Value = 10
$ws.Range("H42").Value = 1290

# Moldavia (row 61)
$ws.Range("E61").Value = 6776
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 576

# Marruecos (row 65)
$ws.Range("B65").Value = 13434
$ws.Range("C65").Value = 146
$ws.Range("D65").Value = 9214
$ws.Range("E65").Value = 3988
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 232

# Consejo Danes para los Refugiados (row 77)
$ws.Range("B77").Value = 7379
$ws.Range("C77").Value = 68
$ws.Range("D77").Value = 2961
$ws.Range("E77").Value = 4236
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 182

# Finlandia (row 78)
$ws.Range("B78").Value = 7248
$ws.Range("C78").Value = 6
$ws.Range("E78").Value = 219

# Hungria (row 96)
$ws.Range("E96").Value = 833
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 589

# Estado de Palestina (row 98)
$ws.Range("E98").Value = 3124
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 12

# Estonia (row 114)
$ws.Range("B114").Value = 1993
$ws.Range("C114").Value = 2
$ws.Range("D114").Value = 1870
$ws.Range("E114").Value = 54

# Eslovenia (row 120)
$ws.Range("B120").Value = 1679
$ws.Range("C120").Value = 29
$ws.Range("E120").Value = 184
